$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31 -- this shifts the existing rows 31..47
# down to 32..48, preserving all of their data/formatting untouched.
$ws.Rows("31:31").Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = "Vega Modelo de Temuco"
$ws.Range("C31").Value = "La Araucanía"
$ws.Range("D31").Value = 44837
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 100112042
$ws.Range("G31").Value = "Locoto"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 50
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = 2500
$ws.Range("N31").Value = "$/kilo"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 2500
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"
